$wb = $excel.ActiveWorkbook

# --- start_price sheet ---
$wsStart = $wb.Worksheets.Item("start_price")
$wsStart.Range("A2").Value = 1.6426

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.00002697150114347342
$wsLinear.Range("B3").Value = -0.1191214929773638
$wsLinear.Range("B4").Value = 0.001721328424401259
$wsLinear.Range("B5").Value = "[1.0, 0.15214199063614664, 0.0024567746449741306, -0.05354909814529638, -0.07315228782721622, -0.016372744574715402, 0.15372033099980228, 0.30059508547459485, 0.12548394789733197, -0.02969963347753591, -0.09714544857305181, -0.07018405305991973, -0.051899849740815164, 0.15227050616451282, 0.30034768904030523, 0.1360158520957131, -0.030272469010306467, -0.10188644862007384, -0.08816246104797394, -0.04440799880289943]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.9513091410197519
$wsNonLinear.Range("B4").Value = 0.0002938029793293068
$wsNonLinear.Range("B5").Value = -0.05622652803176666
$wsNonLinear.Range("B6").Value = 0.001767962173888043
$wsNonLinear.Range("B7").Value = 0.003155959689900679
$wsNonLinear.Range("B8").Value = -0.326588889523496
$wsNonLinear.Range("B9").Value = 0.001672524500372691
$wsNonLinear.Range("B10").Value = "[1.0, 0.1544674970166712, 0.004559652370383741, -0.049658223161043714, -0.06912541428702318, -0.0184221376826904, 0.1516922152660957, 0.29904641062185816, 0.1223992717484115, -0.02969442462103319, -0.0941566456745915, -0.06779237347259578, -0.052760395509828, 0.15068692223492558, 0.297977045128098, 0.1345681922024824, -0.029657986747335035, -0.09985994708396567, -0.08441942880303362, -0.04410857363992214]"
